$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 359.4
$ws.Range("I4").Value = 359.4
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 359.4
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -245.4

$ws.Range("N18").ClearContents()
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -716

$ws.Range("H43").Value = 3163.3333
$ws.Range("I43").Value = 5460
$ws.Range("J43").Value = 866.6667
$ws.Range("K43").Value = 5460
$ws.Range("L43").Value = 866.6667
$ws.Range("M43").Value = -5391
$ws.Range("N43").Value = -1004.6667

$ws.Range("H58").Value = 2215.5264
$ws.Range("I58").Value = 163.18182
$ws.Range("J58").Value = 5037.5
$ws.Range("K58").Value = 489.5454599999999
$ws.Range("L58").Value = 15112.5
$ws.Range("M58").Value = -339.5454599999999
$ws.Range("N58").Value = -15412.5

$ws.Range("H69").Value = 5543.3335
$ws.Range("I69").Value = 3940
$ws.Range("J69").Value = 5743.75
$ws.Range("K69").Value = 11820
$ws.Range("L69").Value = 17231.25
$ws.Range("M69").Value = -10946
$ws.Range("N69").Value = -18979.25

$ws.Range("H72").Value = 5543.3335
$ws.Range("I72").Value = 3940
$ws.Range("J72").Value = 5743.75
$ws.Range("K72").Value = 35460
$ws.Range("L72").Value = 51693.75
$ws.Range("M72").Value = -31092
$ws.Range("N72").Value = -60429.75

$ws.Range("H74").Value = 10264.667
$ws.Range("I74").Value = 14385.556
$ws.Range("J74").Value = 4083.3333
$ws.Range("K74").Value = 14385.556
$ws.Range("L74").Value = 4083.3333
$ws.Range("M74").Value = -13449.556

$ws.Range("H77").Value = 10264.667
$ws.Range("I77").Value = 14385.556
$ws.Range("J77").Value = 4083.3333
$ws.Range("K77").Value = 71927.78
$ws.Range("L77").Value = 20416.6665
$ws.Range("M77").Value = -67247.78

$ws.Range("H100").Value = 100001600
$ws.Range("I100").Value = 1300
$ws.Range("J100").Value = 166668460
$ws.Range("K100").Value = 1300
$ws.Range("L100").Value = 166668460
$ws.Range("M100").Value = -759
$ws.Range("N100").Value = -166669542

$ws.Range("H137").Value = 4286894
$ws.Range("I137").Value = 1725354.9
$ws.Range("J137").Value = 16667666
$ws.Range("K137").Value = 5176064.699999999
$ws.Range("L137").Value = 50002998
$ws.Range("M137").Value = -5173514.699999999
$ws.Range("N137").Value = -50008098

$ws.Range("H138").Value = 2621.7273
$ws.Range("I138").Value = 2294.2856
$ws.Range("J138").Value = 3194.75
$ws.Range("K138").Value = 6882.8568
$ws.Range("L138").Value = 9584.25
$ws.Range("M138").Value = -1742.8568
$ws.Range("N138").Value = -19864.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 30029.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 30029.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 30029.5
$ws.Range("N7").Value = -30257.5

$ws.Range("H61").Value = 1156.9
$ws.Range("I61").Value = 1220
$ws.Range("J61").Value = 799.3333
$ws.Range("K61").Value = 1220
$ws.Range("L61").Value = 799.3333
$ws.Range("M61").Value = -1008
$ws.Range("N61").Value = -1223.3333

$ws.Range("H74").Value = 1048.7838
$ws.Range("I74").Value = 1043.0358
$ws.Range("J74").Value = 1066.6666
$ws.Range("K74").Value = 1043.0358
$ws.Range("L74").Value = 1066.6666
$ws.Range("M74").Value = -169.0358000000001
$ws.Range("N74").Value = -2814.6666

$ws.Range("H77").Value = 1048.7838
$ws.Range("I77").Value = 1043.0358
$ws.Range("J77").Value = 1066.6666
$ws.Range("K77").Value = 5215.179
$ws.Range("L77").Value = 5333.333000000001
$ws.Range("M77").Value = -847.1790000000001
$ws.Range("N77").Value = -14069.333

$ws.Range("H132").Value = 210380.67
$ws.Range("I132").Value = 239815.42
$ws.Range("J132").Value = 4337.3335
$ws.Range("K132").Value = 719446.26
$ws.Range("L132").Value = 13012.0005
$ws.Range("M132").Value = -716916.26
$ws.Range("N132").Value = -18072.0005

$ws.Range("H136").Value = 1156.9
$ws.Range("I136").Value = 1220
$ws.Range("J136").Value = 799.3333
$ws.Range("K136").Value = 3660
$ws.Range("L136").Value = 2397.9999
$ws.Range("M136").Value = -1110
$ws.Range("N136").Value = -7497.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 275172.8
$ws.Range("I134").Value = 430614.56
$ws.Range("J134").Value = 3149.75
$ws.Range("K134").Value = 1291843.68
$ws.Range("L134").Value = 9449.25
$ws.Range("M134").Value = -1289308.68
$ws.Range("N134").Value = -14519.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12822380
$ws.Range("I31").Value = 1525.5294
$ws.Range("J31").Value = 37039548
$ws.Range("K31").Value = 1525.5294
$ws.Range("L31").Value = 37039548
$ws.Range("M31").Value = -1230.5294
$ws.Range("N31").Value = -37040138

$ws.Range("H34").Value = 12822380
$ws.Range("I34").Value = 1525.5294
$ws.Range("J34").Value = 37039548
$ws.Range("K34").Value = 1525.5294
$ws.Range("L34").Value = 37039548
$ws.Range("M34").Value = -1323.5294
$ws.Range("N34").Value = -37039952

$ws.Range("H132").Value = 3558.04
$ws.Range("I132").Value = 2759.6191
$ws.Range("J132").Value = 7749.75
$ws.Range("K132").Value = 8278.8573
$ws.Range("L132").Value = 23249.25
$ws.Range("M132").Value = -5748.8573

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N92").ClearContents()
$ws.Range("H92").Value = 451
$ws.Range("I92").Value = 451
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1353
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -105

$ws.Range("H100").Value = 10108594
$ws.Range("I100").Value = 6000
$ws.Range("J100").Value = 10760374
$ws.Range("K100").Value = 18000
$ws.Range("L100").Value = 32281122
$ws.Range("M100").Value = -17189
$ws.Range("N100").Value = -32282744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2645
$ws.Range("I132").Value = 1876
$ws.Range("J132").Value = 3926.6667
$ws.Range("K132").Value = 5628
$ws.Range("L132").Value = 11780.0001
$ws.Range("M132").Value = -3098
$ws.Range("N132").Value = -16840.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 25644.4
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 25644.4
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 25644.4
$ws.Range("N95").Value = -31136.4

$ws.Range("H132").Value = 2808.3044
$ws.Range("I132").Value = 2545.2307
$ws.Range("J132").Value = 3150.3
$ws.Range("K132").Value = 7635.6921
$ws.Range("L132").Value = 9450.900000000001
$ws.Range("M132").Value = -5105.6921
$ws.Range("N132").Value = -14510.9

$ws.Range("H136").Value = 1545.027
$ws.Range("I136").Value = 1376.1714
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 4128.5142
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -1578.5142
$ws.Range("N136").Value = -18600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1784.0416
$ws.Range("I126").Value = 1734
$ws.Range("J126").Value = 1934.1666
$ws.Range("K126").Value = 5202
$ws.Range("L126").Value = 5802.4998
$ws.Range("M126").Value = -2732
$ws.Range("N126").Value = -10742.4998

$ws.Range("H132").Value = 2202.606
$ws.Range("I132").Value = 2253.0435
$ws.Range("J132").Value = 2086.6
$ws.Range("K132").Value = 6759.130500000001
$ws.Range("L132").Value = 6259.799999999999
$ws.Range("M132").Value = -4229.130500000001

$ws.Range("H136").Value = 1320.8704
$ws.Range("I136").Value = 1297.84
$ws.Range("J136").Value = 1608.75
$ws.Range("K136").Value = 3893.52
$ws.Range("L136").Value = 4826.25
$ws.Range("M136").Value = -1343.52
$ws.Range("N136").Value = -9926.25
